$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New query text for the existing "CasesTab" row (row 2, columns B & C)
# ---------------------------------------------------------------------------
$casesQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "UNKNOWN"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE c.race = "UNKNOWN"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

$filesQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
  WHERE c.race = "UNKNOWN"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# ---------------------------------------------------------------------------
# Add a new row 3 for the "FilesTab" (write the label first so the new
# shared-string table ends up ordered the same way as the source workbook:
# CasesTab, FilesTab, new-cases-query, new-stat-query, new-files-query)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "FilesTab"

# Update row 2 (CasesTab) with the new query text
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# Row 2 grows a bit taller to fit the new query text
$ws.Rows.Item(2).RowHeight = 195

$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC05_Trials_Filter_Race-Unknown_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC05_Trials_Filter_Race-Unknown_WebData.xlsx"

$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Row 3 needs the maximum row height to display the long query
$ws.Rows.Item(3).RowHeight = 409.5

# ---------------------------------------------------------------------------
# View-level tweaks: zoom to 70% and move the active selection to C3
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 70
$ws.Range("C3").Select()
